$d = $word.ActiveDocument

# --- Change 1: "country Survey" (Heading1) -> "country"
# The trailing " Survey" text (spread across 3 runs) is removed; the
# following Heading2 paragraph ("Sample Description and Weighting
# Procedures") is left untouched.
$d.Content.Find.Execute("country Survey", $false, $false, $false, $false, $false, $true, 1, $false, "country", 2) | Out-Null

# --- Change 2: VML straight-connector line - refresh z-index / drop the
# top+bottom wrap-distance properties from its `style` attribute. The
# shape is legacy VML (w:pict), so it is not reachable via the Shapes
# collection - patch the paragraph's raw OOXML instead, changing only the
# `style` attribute and leaving the (large) o:gfxdata fallback blob intact.
$pictXml = @'
<w:p w14:paraId="0C932D5E" w14:textId="412FB6D6" w:rsidR="00284227" w:rsidRDefault="00D77B6B" w:rsidP="00284227"><w:r><w:rPr><w:noProof/></w:rPr><w:pict w14:anchorId="0085882E"><v:line id="Straight Connector 1" o:spid="_x0000_s2050" style="position:absolute;z-index:251657728;visibility:visible;mso-wrap-style:square;mso-width-percent:0;mso-height-percent:0;mso-wrap-distance-left:9pt;mso-wrap-distance-right:9pt;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;mso-width-percent:0;mso-height-percent:0;mso-width-relative:margin;mso-height-relative:margin" from="-1.6pt,5.45pt" to="453.85pt,5.45pt" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#xD;&#xA;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#xD;&#xA;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#xD;&#xA;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#xD;&#xA;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#xD;&#xA;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#xD;&#xA;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#xD;&#xA;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#xD;&#xA;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#xD;&#xA;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#xD;&#xA;IQBi/0fCxwEAAHgDAAAOAAAAZHJzL2Uyb0RvYy54bWysU8tu2zAQvBfoPxC815Jdu3EFyznYSC9B&#xD;&#xA;GyDpB6wpUiLKF7isJf99l/QjQXorqgNB7nJnd4ajzf1kDTvKiNq7ls9nNWfSCd9p17f858vDpzVn&#xD;&#xA;mMB1YLyTLT9J5Pfbjx82Y2jkwg/edDIyAnHYjKHlQ0qhqSoUg7SAMx+ko6Ty0UKiY+yrLsJI6NZU&#xD;&#xA;i7r+Uo0+diF6IREpuj8n+bbgKyVF+qEUysRMy2m2VNZY1kNeq+0Gmj5CGLS4jAH/MIUF7ajpDWoP&#xD;&#xA;CdjvqP+CslpEj16lmfC28kppIQsHYjOv37F5HiDIwoXEwXCTCf8frPh+3LmnmEcXk3sOj178Qub8&#xD;&#xA;bgDXyzLAyynQw82zVNUYsLmV5AOGc/Gkos0gxIhNRd7TTV45JSYouLpbLxfzFWfimquguRaGiOmb&#xD;&#xA;9JblTcuNdpk5NHB8xJRbQ3O9ksPOP2hjyusZx0aa7/Nq/ZWggUykDCTa2tC1HF3PGZie3ClSLJDo&#xD;&#xA;je5yeQbC2B92JrIjkEOWy7vFblmYUubttdx7Dzic75XU2TtWJzKw0bbl6zp/l2rjMrosFrwweNUr&#xD;&#xA;7w6+Oz3Fq6j0vIXjxYrZP2/PRfrXH2b7BwAA//8DAFBLAwQUAAYACAAAACEAxgY3E+EAAAANAQAA&#xD;&#xA;DwAAAGRycy9kb3ducmV2LnhtbExPTU/DMAy9I/EfIiNx25IV1LKu6QRDCKSdNhBcsyY01RqnJNlW&#xD;&#xA;+PUYcYCLJb9nv49qObqeHU2InUcJs6kAZrDxusNWwsvzw+QGWEwKteo9GgmfJsKyPj+rVKn9CTfm&#xD;&#xA;uE0tIxGMpZJgUxpKzmNjjVNx6geDxL374FSiNbRcB3UicdfzTIicO9UhOVg1mJU1zX57cBJebZGt&#xD;&#xA;v/L93WNQm/z6afXRzN5yKS8vxvsFjdsFsGTG9PcBPx0oP9QUbOcPqCPrJUyuMrokXMyBET8XRQFs&#xD;&#xA;9wvwuuL/W9TfAAAA//8DAFBLAQItABQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAAAAAAAAAAAAA&#xD;&#xA;AAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhADj9If/WAAAAlAEAAAsAAAAA&#xD;&#xA;AAAAAAAAAAAALwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhAGL/R8LHAQAAeAMAAA4AAAAA&#xD;&#xA;AAAAAAAAAAAALgIAAGRycy9lMm9Eb2MueG1sUEsBAi0AFAAGAAgAAAAhAMYGNxPhAAAADQEAAA8A&#xD;&#xA;AAAAAAAAAAAAAAAAIQQAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAABAAEAPMAAAAvBQAAAAA=&#xD;&#xA;" strokecolor="#4472c4" strokeweight="1.07pt"><v:stroke joinstyle="miter"/></v:line></w:pict></w:r></w:p>
'@
$linePara = $d.Paragraphs(3)
$linePara.Range.InsertXML($pictXml)

# --- Change 3: merge " with " + "the majority of" + " students in" (and
# drop the two w:proofErr markers) into a single run reading
# " with the majority of students in". The preceding run ("All classes")
# shares identical run formatting and directly abuts this text, so a
# straight Find/Replace over the whole phrase would also swallow "All
# classes" into the new run. Temporarily drop a bookmark right at the
# boundary to keep the two runs from being coalesced, then remove it
# once the replacement is done.
$boundary = $d.Content
$boundary.Find.ClearFormatting()
$boundary.Find.Execute("All classes", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sepRange = $d.Range($boundary.End, $boundary.End)
$sepRange.Bookmarks.Add("zzTempSeparator") | Out-Null

$target = $d.Content
$target.Find.ClearFormatting()
$target.Find.Execute(" with the majority of students in", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target = $d.Range($target.Start, $target.End)
$target.Find.Execute(" with the majority of students in", $false, $false, $false, $false, $false, $true, 1, $false, " with the majority of students in", 2) | Out-Null

if ($d.Bookmarks.Exists("zzTempSeparator")) {
    $d.Bookmarks("zzTempSeparator").Delete()
}

Write-Output "done"
